$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.929.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.551.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.60%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.05'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.486'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.249'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.61'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0588'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.772.00'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.552.68'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.47%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.04%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.931.93'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.81'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.27'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0688'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.02'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.12'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.88'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.47%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.87'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0464'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.06%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.405.70'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.78%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.955'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.54%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.520'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.989'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.44%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.91%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.37%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.81%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.686.13'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.23'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.50%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0953'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.14%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₇0966'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.14%  '
